$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Accuracy of correct prediction side"
$ws.Range("B12").Value = 0.493
